$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before the current row 135, shifting existing
# rows 135-204 down to 137-206 (dimension grows from T204 to T206).
$ws.Rows.Item(135).EntireRow.Insert()
$ws.Rows.Item(135).EntireRow.Insert()

# Fill in the two newly inserted rows (135 and 136) with the new weekly data.
$ws.Cells.Item(135, 1).Value = 10
$ws.Cells.Item(135, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(135, 3).Value = "La Araucanía"
$ws.Cells.Item(135, 4).Value = 44455
$ws.Cells.Item(135, 5).Value = 9
$ws.Cells.Item(135, 6).Value = "Fruta"
$ws.Cells.Item(135, 7).Value = 100108
$ws.Cells.Item(135, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(135, 9).Value = 100108002
$ws.Cells.Item(135, 10).Value = "Mango"
$ws.Cells.Item(135, 11).Value = "Sin especificar"
$ws.Cells.Item(135, 12).Value = "Primera"
$ws.Cells.Item(135, 13).Value = 55
$ws.Cells.Item(135, 14).Value = 9000
$ws.Cells.Item(135, 15).Value = 9000
$ws.Cells.Item(135, 16).Value = 9000
$ws.Cells.Item(135, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(135, 18).Value = "Brasil"
$ws.Cells.Item(135, 19).Value = 2250
$ws.Cells.Item(135, 20).Value = 4

$ws.Cells.Item(136, 1).Value = 10
$ws.Cells.Item(136, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(136, 3).Value = "La Araucanía"
$ws.Cells.Item(136, 4).Value = 44455
$ws.Cells.Item(136, 5).Value = 9
$ws.Cells.Item(136, 6).Value = "Fruta"
$ws.Cells.Item(136, 7).Value = 100108
$ws.Cells.Item(136, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(136, 9).Value = 100108002
$ws.Cells.Item(136, 10).Value = "Mango"
$ws.Cells.Item(136, 11).Value = "Sin especificar"
$ws.Cells.Item(136, 12).Value = "Segunda"
$ws.Cells.Item(136, 13).Value = 45
$ws.Cells.Item(136, 14).Value = 6000
$ws.Cells.Item(136, 15).Value = 6000
$ws.Cells.Item(136, 16).Value = 6000
$ws.Cells.Item(136, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(136, 18).Value = "Brasil"
$ws.Cells.Item(136, 19).Value = 1500
$ws.Cells.Item(136, 20).Value = 4

# Apply the same date number format (style) used by the other date cells
# in column D to the two new date cells.
$ws.Cells.Item(135, 4).NumberFormat = $ws.Cells.Item(137, 4).NumberFormat
$ws.Cells.Item(136, 4).NumberFormat = $ws.Cells.Item(137, 4).NumberFormat
